# Sync attendance_reports: reorder "Recorded By" name lists in column G
# so that "System" appears before the user email(s) it was recorded alongside.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value   = "system, System, backup@backdoor.com"
$ws.Range("G3").Value   = "System, dnasr281@gmail.com"
$ws.Range("G5").Value   = "System, backup@backdoor.com"
$ws.Range("G6").Value   = "System, dnasr281@gmail.com"
$ws.Range("G7").Value   = "System, admin@admin.com"
$ws.Range("G8").Value   = "System, backup@backdoor.com"

$ws.Range("G28").Value  = "system, System, backup@backdoor.com"
$ws.Range("G29").Value  = "System, dnasr281@gmail.com"
$ws.Range("G31").Value  = "System, backup@backdoor.com"
$ws.Range("G32").Value  = "System, dnasr281@gmail.com"
$ws.Range("G33").Value  = "System, admin@admin.com"
$ws.Range("G34").Value  = "System, backup@backdoor.com"

$ws.Range("G54").Value  = "system, System, backup@backdoor.com"
$ws.Range("G55").Value  = "System, dnasr281@gmail.com"
$ws.Range("G57").Value  = "System, backup@backdoor.com"
$ws.Range("G58").Value  = "System, dnasr281@gmail.com"
$ws.Range("G59").Value  = "System, admin@admin.com"
$ws.Range("G60").Value  = "System, backup@backdoor.com"

$ws.Range("G80").Value  = "System, backup@backdoor.com"
$ws.Range("G81").Value  = "System, backup@backdoor.com"
$ws.Range("G82").Value  = "System, backup@backdoor.com"
$ws.Range("G87").Value  = "dnasr281@gmail.com, admin@admin.com"

$ws.Range("G106").Value = "System, backup@backdoor.com"
$ws.Range("G107").Value = "System, backup@backdoor.com"
$ws.Range("G108").Value = "System, backup@backdoor.com"
$ws.Range("G113").Value = "dnasr281@gmail.com, admin@admin.com"

$ws.Range("G132").Value = "System, backup@backdoor.com"
$ws.Range("G133").Value = "System, backup@backdoor.com"
$ws.Range("G134").Value = "System, backup@backdoor.com"
$ws.Range("G139").Value = "dnasr281@gmail.com, admin@admin.com"
